$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (tab) from "vaccinatietoestand" to "PNEU"
$ws.Name = "PNEU"

# 2. Update the four cells that held "primaire serie" to the new quoted
#    capitalised literal "Primaire serie" (matches style of other DMN literals
#    like "Pneu", "Volledig afgesloten", etc.)
$ws.Range("I10").Value = '"Primaire serie"'
$ws.Range("I11").Value = '"Primaire serie"'
$ws.Range("I15").Value = '"Primaire serie"'
$ws.Range("I17").Value = '"Primaire serie"'

# 3. Remove the data validations that were on column G
$ws.Cells.Validation.Delete()

# 4. Collapse the AutoFilter down to just the header row (A1:I1)
$ws.AutoFilterMode = $false
$ws.Range("A1:I1").AutoFilter()

# 5. Point the hidden _FilterDatabase defined name at the new, smaller range
$fdb = $wb.Names.Item("PNEU!_FilterDatabase")
$fdb.RefersTo = "=PNEU!`$A`$1:`$I`$1"

# 6. Select the whole of row 1, mirroring the author's last interaction
$ws.Rows(1).Select()
